$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.320.56"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "1.874.80"
$ws.Range("E3").Value = "  +0.80%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7121"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.32"
$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3109"
$ws.Range("E8").Value = "  +1.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07744"
$ws.Range("E9").Value = "  +0.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.00"
$ws.Range("E10").Value = "  +0.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08471"
$ws.Range("E11").Value = "  +2.61%  "

$ws.Range("D12").Value = "1.891.95"
$ws.Range("E12").Value = "  +0.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.203"
$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7106"
$ws.Range("E14").Value = "  -0.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.29"
$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("D16").Value = "29.329.89"
$ws.Range("E16").Value = "  +0.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008279"
$ws.Range("E17").Value = "  +6.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.999"
$ws.Range("E18").Value = "  +2.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.26"
$ws.Range("E19").Value = "  -0.50%  "

$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.123.45"
$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.20"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.839"
$ws.Range("E23").Value = "  -1.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1611"
$ws.Range("E25").Value = "  +2.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.69"
$ws.Range("E26").Value = "  +0.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.012"
$ws.Range("E27").Value = "  +1.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.47"
$ws.Range("E28").Value = "  +1.18%  "

$ws.Range("E29").Value = "  +1.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.405"
$ws.Range("E30").Value = "  +1.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.346"
$ws.Range("E31").Value = "  +6.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.269"
$ws.Range("E32").Value = "  -3.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05264"
$ws.Range("E33").Value = "  +1.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.929"
$ws.Range("E34").Value = "  +1.34%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7448"
$ws.Range("E36").Value = "  +2.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.682"
$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("E38").Value = "  +1.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.721"

$ws.Range("D40").Value = "1.172.62"
$ws.Range("E40").Value = "  +2.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.375"
$ws.Range("E41").Value = "  +4.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.06"
$ws.Range("E42").Value = "  +1.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8867"
$ws.Range("E43").Value = "  -1.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.99"
$ws.Range("E44").Value = "  +5.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9995"
$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("D46").Value = "2.021.23"
$ws.Range("E46").Value = "  +0.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.816"
$ws.Range("E47").Value = "  +3.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5205"
$ws.Range("E48").Value = "  -0.64%  "

$ws.Range("E49").Value = "  +1.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.383"
$ws.Range("E50").Value = "  +1.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4302"
$ws.Range("E51").Value = "  +1.28%  "
